# Update NIFTY 50 losers indices worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new column header "To Year High Profit" before "Total Year Up" ---
# Current header H1 = "Total Year Up". We need H1 = "To Year High Profit" and I1 = "Total Year Up"
# Copy H1's formatting (bold/border/centered) to I1, then set the new values/labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "Total Year Up"
$ws.Range("H1").Value = "To Year High Profit"

# --- Full target data (rows 2-14), columns A-I ---
# A ISIN, B Symbol, C Industry, D Current, E Year High, F Year Low,
# G Total Year Down, H To Year High Profit, I Total Year Up
$data = @(
    @("INE095A01012","INDUSINDBK","Private Sector Bank",996,1694.5,926.45,41.22159929182649,70.13052208835342,7.507150952560848),
    @("INE155A01022","TATAMOTORS","Passenger Cars & Utility Vehicles",789,1179,717.7,33.07888040712468,49.42965779467681,9.934513027727455),
    @("INE158A01026","HEROMOTOCO","2/3 Wheelers",4240.05,6246.25,3929.85,32.11847108264959,47.31547977028572,7.893431046986543),
    @("INE423A01024","ADANIENT","Trading - Minerals",2562,3743.9,2025,31.56868506103262,46.13192818110852,26.51851851851852),
    @("INE021A01026","ASIANPAINT","Paints",2333.8,3401.25,2256.5,31.3840499816244,45.73870940097693,3.425659206736098),
    @("INE917I01010","BAJAJ-AUTO","2/3 Wheelers",8965,12774,6604,29.81838108658213,42.48745119910764,35.75105996365839),
    @("INE522F01014","COALINDIA","Coal",394.1,543.55,368,27.49517063747584,37.92184724689163,7.092391304347823),
    @("INE742F01042","ADANIPORTS","Port & Port services",1197.8,1621.4,995.65,26.12557049463427,35.36483553180834,20.3033194395621),
    @("INE192A01025","TATACONSUM","Tea & Coffee",936,1256.44,882.9,25.50380439973258,34.23504273504274,6.014271151885842),
    @("INE216A01030","BRITANNIA","Packaged Foods",4839.85,6469.9,4641,25.19436158209554,33.67976280256617,4.284636931695762),
    @("INE123W01016","SBILIFE","Life Insurance",1450,1936,1307.7,25.10330578512396,33.51724137931033,10.88170069587826),
    @("INE081A01020","TATASTEEL","Iron & Steel",138.33,184.6,128.2,25.06500541711808,33.44899877105472,7.901716068642761),
    @("INE213A01029","ONGC","Oil Exploration & Production",258.65,345,204.5,25.02898550724638,33.38488304658807,26.47921760391196)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
